# Fixzone.xlsx update — 2019-05-10 kl. 17:38
# Adds a "Date Friday" calendar block (columns J/K) next to the weekly
# tracker, bumps week-19's Friday (H8/I8) to "done", and refreshes the
# running totals that depend on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) New text labels -----------------------------------------------
# Written in this exact order so the shared-string table fills up the
# same way the original author's keystrokes would have (Maj..September,
# then the "Date Friday" header typed in last).
$ws.Range("J8").Value  = "Maj"
$ws.Range("J12").Value = "Juni"
$ws.Range("J16").Value = "Juli"
$ws.Range("J20").Value = "Augisti"
$ws.Range("J25").Value = "September"
$ws.Range("K7").Value  = "Date Friday"

# --- 2) Friday-date numbers in column K (rows 8-28) --------------------
$fridayDates = @(10,17,24,31,7,14,21,28,5,12,19,26,1,8,15,22,29,6,13,20,27)
for ($i = 0; $i -lt $fridayDates.Length; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 11).Value = $fridayDates[$i]
}

# --- 3) Formatting ------------------------------------------------------
# Month labels (incl. the blank filler cells in the same column range)
# are bold + centred.
$monthRange = $ws.Range("J8:J25")
$monthRange.Font.Bold = $true
$monthRange.HorizontalAlignment = -4108   # xlCenter

# Bottom three filler cells (J26:J28) are centred but NOT bold.
$ws.Range("J26:J28").HorizontalAlignment = -4108

# "Date Friday" header is bold (no centring).
$ws.Range("K7").Font.Bold = $true

# Friday-date numbers are centred (not bold) — matches the rest of col K.
$ws.Range("K8:K28").HorizontalAlignment = -4108

# Column J needs to be wide enough for "September" / "Date Friday".
$ws.Columns.Item(10).ColumnWidth = 12.3

# --- 4) Week 19 (row 8): Friday got done -------------------------------
# Copy the "done" look from the Friday cell one row up (E4 carries the
# same green-fill / thin-border style used across the tracker) and mark
# the day count.
$ws.Range("E4").Copy()
$ws.Range("H8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I8").Value = 4

# --- 5) Selection, for parity with the saved workbook state ------------
[void]$ws.Range("L8").Select()
